$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated 2024 year-end standings (rows 191-200), columns:
# A=Yr B=Person C=SRank D=Points E=PointsBonus F=Bonus G=Chips H=Winnings I=Takehome J=PersStatus K=pers_personid

$rows = @(
    @{ Row = 191; Person = "Richard";  SRank = 1;  Points = 63; PointsBonus = 0; Bonus = 63; Chips = 206300; Winnings = 200; Takehome = 100;  PersonId = 366 },
    @{ Row = 192; Person = "Mark";     SRank = 2;  Points = 47; PointsBonus = 0; Bonus = 47; Chips = 140750; Winnings = 180; Takehome = 70;   PersonId = 361 },
    @{ Row = 193; Person = "Andy";     SRank = 3;  Points = 40; PointsBonus = 0; Bonus = 40; Chips = 138000; Winnings = 130; Takehome = 30;   PersonId = 349 },
    @{ Row = 194; Person = "Anthony";  SRank = 4;  Points = 32; PointsBonus = 0; Bonus = 32; Chips = 108750; Winnings = 40;  Takehome = -40;  PersonId = 350 },
    @{ Row = 195; Person = "Matt";     SRank = 5;  Points = 30; PointsBonus = 0; Bonus = 30; Chips = 116550; Winnings = 40;  Takehome = -70;  PersonId = 362 },
    @{ Row = 196; Person = "Prashant"; SRank = 6;  Points = 28; PointsBonus = 0; Bonus = 28; Chips = 97650;  Winnings = 130; Takehome = 40;   PersonId = 365 },
    @{ Row = 197; Person = "Pepe";     SRank = 7;  Points = 24; PointsBonus = 0; Bonus = 24; Chips = 77350;  Winnings = 40;  Takehome = -40;  PersonId = 364 },
    @{ Row = 198; Person = "Jon";      SRank = 8;  Points = 23; PointsBonus = 0; Bonus = 23; Chips = 76650;  Winnings = 20;  Takehome = -90;  PersonId = 357 },
    @{ Row = 199; Person = "Maisy";    SRank = 9;  Points = 20; PointsBonus = 0; Bonus = 20; Chips = 74950;  Winnings = 50;  Takehome = -40;  PersonId = 360 },
    @{ Row = 200; Person = "Alex";     SRank = 10; Points = 19; PointsBonus = 1; Bonus = 20; Chips = 69950;  Winnings = 80;  Takehome = 40;   PersonId = 348 }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 2).Value = $r.Person
    $ws.Cells.Item($r.Row, 3).Value = $r.SRank
    $ws.Cells.Item($r.Row, 4).Value = $r.Points
    $ws.Cells.Item($r.Row, 5).Value = $r.PointsBonus
    $ws.Cells.Item($r.Row, 6).Value = $r.Bonus
    $ws.Cells.Item($r.Row, 7).Value = $r.Chips
    $ws.Cells.Item($r.Row, 8).Value = $r.Winnings
    $ws.Cells.Item($r.Row, 9).Value = $r.Takehome
    $ws.Cells.Item($r.Row, 11).Value = $r.PersonId
}
